$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object 'object[,]' 1,20
$row2[0,0] = "ECs"
$row2[0,1] = "Gal"
$row2[0,2] = "Galr1"
$row2[0,3] = "sCs"
$row2[0,4] = 3
$row2[0,5] = 1
$row2[0,6] = 4.926742666666667
$row2[0,7] = 14.780228
$row2[0,8] = 0.666446655488151
$row2[0,9] = 0.6664466554881509
$row2[0,10] = 3
$row2[0,11] = 1
$row2[0,12] = 0.8078946666666668
$row2[0,13] = 2.423684
$row2[0,14] = 1
$row2[0,15] = 1
$row2[0,16] = 3.980289124439111
$row2[0,17] = 35.822602119952
$row2[0,18] = 0.666446655488151
$row2[0,19] = 0.6664466554881509
$ws.Range("A2:T2").Value = $row2

$row3 = New-Object 'object[,]' 1,20
$row3[0,0] = "FAPs"
$row3[0,1] = "Gal"
$row3[0,2] = "Galr1"
$row3[0,3] = "sCs"
$row3[0,4] = 2
$row3[0,5] = 0.6666666666666666
$row3[0,6] = 0.48067
$row3[0,7] = 1.44201
$row3[0,8] = 0.06502083335118163
$row3[0,9] = 0.06502083335118163
$row3[0,10] = 3
$row3[0,11] = 1
$row3[0,12] = 0.8078946666666668
$row3[0,13] = 2.423684
$row3[0,14] = 1
$row3[0,15] = 1
$row3[0,16] = 0.3883307294266667
$row3[0,17] = 3.49497656484
$row3[0,18] = 0.06502083335118163
$row3[0,19] = 0.06502083335118163
$ws.Range("A3:T3").Value = $row3

$row4 = New-Object 'object[,]' 1,20
$row4[0,0] = "sCs"
$row4[0,1] = "Gal"
$row4[0,2] = "Galr1"
$row4[0,3] = "sCs"
$row4[0,4] = 3
$row4[0,5] = 1
$row4[0,6] = 1.985141
$row4[0,7] = 5.955423
$row4[0,8] = 0.2685325111606675
$row4[0,9] = 0.2685325111606675
$row4[0,10] = 3
$row4[0,11] = 1
$row4[0,12] = 0.8078946666666668
$row4[0,13] = 2.423684
$row4[0,14] = 1
$row4[0,15] = 1
$row4[0,16] = 1.603784826481333
$row4[0,17] = 14.434063438332
$row4[0,18] = 0.2685325111606675
$row4[0,19] = 0.2685325111606675
$ws.Range("A4:T4").Value = $row4
